$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.425.35"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").Value = "3.078.20"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.11"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.27"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.074.99"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.25"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("E12").Value = "  -2.56%  "

$ws.Range("E13").Value = "  -3.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.76"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("E15").Value = "  -1.45%  "

$ws.Range("D16").Value = "3.591.56"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").Value = "66.374.56"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("E18").Value = "  -2.68%  "

$ws.Range("D19").Value = "3.077.84"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.57"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.96"
$ws.Range("E21").Value = "  +2.67%  "

$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("E23").Value = "  -3.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.48"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("E25").Value = "  -4.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("E30").Value = "  -5.00%  "

$ws.Range("E31").Value = "  -3.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.76"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("E33").Value = "  -3.64%  "

$ws.Range("D34").Value = "0.0₃0909"
$ws.Range("E34").Value = "  -4.33%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.07"
$ws.Range("E36").Value = "  +2.07%  "

$ws.Range("E37").Value = "  -4.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.941"
$ws.Range("E38").Value = "  -3.32%  "

$ws.Range("E39").Value = "  -1.24%  "

$ws.Range("E40").Value = "  -3.41%  "

$ws.Range("E41").Value = "  -4.60%  "

$ws.Range("E42").Value = "  -4.38%  "

$ws.Range("D43").Value = "2.772.14"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("E45").Value = "  -2.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.66"
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "365.18"
$ws.Range("E47").Value = "  -4.45%  "

$ws.Range("E49").Value = "  -2.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("E51").Value = "  -2.03%  "
